$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02879502170626
$ws.Range("D2").Value = 1.033830891926608
$ws.Range("E2").Value = 1.028720230191227
$ws.Range("I2").Value = 1.035571040128027
$ws.Range("J2").Value = 1.033945360437227
$ws.Range("K2").Value = 1.036632111225898
$ws.Range("L2").Value = 1.031536226695657
$ws.Range("N2").Value = 1.035413681444249

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02960165680806
$ws.Range("D3").Value = 1.034440068040579
$ws.Range("E3").Value = 1.02940109902744
$ws.Range("I3").Value = 1.035759800590749
$ws.Range("J3").Value = 1.034393487360341
$ws.Range("K3").Value = 1.037050882515991
$ws.Range("L3").Value = 1.032025433128386
$ws.Range("N3").Value = 1.035862444758984

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030123934658339
$ws.Range("D4").Value = 1.034834381236086
$ws.Range("E4").Value = 1.029842352762316
$ws.Range("I4").Value = 1.035880673923915
$ws.Range("J4").Value = 1.034683121708589
$ws.Range("K4").Value = 1.037321297242682
$ws.Range("L4").Value = 1.032341989910051
$ws.Range("N4").Value = 1.036152490421228

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030343577795149
$ws.Range("D5").Value = 1.035000181099876
$ws.Range("E5").Value = 1.030028018218281
$ws.Range("I5").Value = 1.035931184810866
$ws.Range("J5").Value = 1.03480480277005
$ws.Range("K5").Value = 1.037434844568567
$ws.Range("L5").Value = 1.032475070766454
$ws.Range("N5").Value = 1.036274344283753

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03038046134078
$ws.Range("D6").Value = 1.035028021366974
$ws.Range("E6").Value = 1.030059201733508
$ws.Range("I6").Value = 1.035939647943761
$ws.Range("J6").Value = 1.034825228754535
$ws.Range("K6").Value = 1.037453901706474
$ws.Range("L6").Value = 1.032497415609321
$ws.Range("N6").Value = 1.036294799275479

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030126869238198
$ws.Range("D7").Value = 1.034836596544085
$ws.Range("E7").Value = 1.029844832997505
$ws.Range("I7").Value = 1.035881350049487
$ws.Range("J7").Value = 1.034684747937683
$ws.Range("K7").Value = 1.037322814999885
$ws.Range("L7").Value = 1.032343768143231
$ws.Range("N7").Value = 1.036154118959755

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029067558249217
$ws.Range("D8").Value = 1.034036736678462
$ws.Range("E8").Value = 1.028950189829544
$ws.Range("I8").Value = 1.035635094494137
$ws.Range("J8").Value = 1.034096875337522
$ws.Range("K8").Value = 1.03677375172276
$ws.Range("L8").Value = 1.031701553897513
$ws.Range("N8").Value = 1.035565411513071

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027203527715436
$ws.Range("D9").Value = 1.032628400724347
$ws.Range("E9").Value = 1.027379050850588
$ws.Range("I9").Value = 1.035191492904259
$ws.Range("J9").Value = 1.033058470985553
$ws.Range("K9").Value = 1.035802016462618
$ws.Range("L9").Value = 1.030570010141306
$ws.Range("N9").Value = 1.034525532507816

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025962697697182
$ws.Range("D10").Value = 1.031690372404063
$ws.Range("E10").Value = 1.026335323610398
$ws.Range("I10").Value = 1.034889309799399
$ws.Range("J10").Value = 1.032364600532352
$ws.Range("K10").Value = 1.035151439767282
$ws.Range("L10").Value = 1.029815812088517
$ws.Range("N10").Value = 1.033830676678983

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025425866375414
$ws.Range("D11").Value = 1.031284423455613
$ws.Range("E11").Value = 1.025884278858896
$ws.Range("I11").Value = 1.034756942409548
$ws.Range("J11").Value = 1.032063783748245
$ws.Range("K11").Value = 1.03486909692326
$ws.Range("L11").Value = 1.029489292683794
$ws.Range("N11").Value = 1.033529432700538

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025226533571676
$ws.Range("D12").Value = 1.031133671387357
$ws.Range("E12").Value = 1.025716877268014
$ws.Range("I12").Value = 1.034707547734599
$ws.Range("J12").Value = 1.031951993287356
$ws.Range("K12").Value = 1.03476412760478
$ws.Range("L12").Value = 1.029368018272376
$ws.Range("N12").Value = 1.033417483484372

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025269287928354
$ws.Range("D13").Value = 1.031166006583804
$ws.Range("E13").Value = 1.025752779269672
$ws.Range("I13").Value = 1.034718153342792
$ws.Range("J13").Value = 1.031975975130192
$ws.Range("K13").Value = 1.034786648139575
$ws.Range("L13").Value = 1.029394031595608
$ws.Range("N13").Value = 1.033441499384175

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025409388020715
$ws.Range("D14").Value = 1.031271961508757
$ws.Range("E14").Value = 1.025870438595421
$ws.Range("I14").Value = 1.03475286406701
$ws.Range("J14").Value = 1.032054544201994
$ws.Range("K14").Value = 1.034860422048887
$ws.Range("L14").Value = 1.029479267902458
$ws.Range("N14").Value = 1.033520180033071

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025495717671463
$ws.Range("D15").Value = 1.031337248587679
$ws.Range("E15").Value = 1.025942950535078
$ws.Range("I15").Value = 1.034774220364831
$ws.Range("J15").Value = 1.0321029461221
$ws.Range("K15").Value = 1.034905864092616
$ws.Range("L15").Value = 1.029531786096799
$ws.Range("N15").Value = 1.033568650689456

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025998335182329
$ws.Range("D16").Value = 1.031717318767279
$ws.Range("E16").Value = 1.026365277005201
$ws.Range("I16").Value = 1.034898062632418
$ws.Range("J16").Value = 1.032384557146728
$ws.Range("K16").Value = 1.035170164588784
$ws.Range("L16").Value = 1.029837483351048
$ws.Range("N16").Value = 1.03385066163404

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026313737225066
$ws.Range("D17").Value = 1.031955787907723
$ws.Range("E17").Value = 1.026630432473682
$ws.Range("I17").Value = 1.034975339235995
$ws.Range("J17").Value = 1.032561107165532
$ws.Range("K17").Value = 1.035335783234891
$ws.Range("L17").Value = 1.030029254444497
$ws.Range("N17").Value = 1.034027462374121

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026497749809111
$ws.Range("D18").Value = 1.032094904362763
$ws.Range("E18").Value = 1.026785179499613
$ws.Range("I18").Value = 1.035020266645981
$ws.Range("J18").Value = 1.032664050373649
$ws.Range("K18").Value = 1.035432324058842
$ws.Range("L18").Value = 1.030141116329107
$ws.Range("N18").Value = 1.034130551773402

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026560500740823
$ws.Range("D19").Value = 1.032142343058557
$ws.Range("E19").Value = 1.026837958821339
$ws.Range("I19").Value = 1.035035560818802
$ws.Range("J19").Value = 1.032699145297479
$ws.Range("K19").Value = 1.035465231472911
$ws.Range("L19").Value = 1.030179259154312
$ws.Range("N19").Value = 1.034165696536049

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026279893016831
$ws.Range("D20").Value = 1.031930200194682
$ws.Range("E20").Value = 1.026601974858416
$ws.Range("I20").Value = 1.034967063358157
$ws.Range("J20").Value = 1.032542168674622
$ws.Range("K20").Value = 1.035318020296764
$ws.Range("L20").Value = 1.030008678688431
$ws.Range("N20").Value = 1.034008496988382

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025368130082454
$ws.Range("D21").Value = 1.031240759430762
$ws.Range("E21").Value = 1.025835787072966
$ws.Range("I21").Value = 1.034742648898088
$ws.Range("J21").Value = 1.032031409040803
$ws.Range("K21").Value = 1.034838700072256
$ws.Range("L21").Value = 1.029454167662729
$ws.Range("N21").Value = 1.033497012017299

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024795275834943
$ws.Range("D22").Value = 1.030807486659215
$ws.Range("E22").Value = 1.025354845156333
$ws.Range("I22").Value = 1.034600234604931
$ws.Range("J22").Value = 1.031709964069073
$ws.Range("K22").Value = 1.034536785683485
$ws.Range("L22").Value = 1.029105579854873
$ws.Range("N22").Value = 1.033175110556837

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025098917534555
$ws.Range("D23").Value = 1.031037152604631
$ws.Range("E23").Value = 1.025609725980479
$ws.Range("I23").Value = 1.034675855576664
$ws.Range("J23").Value = 1.031880397086321
$ws.Range("K23").Value = 1.03469688758327
$ws.Range("L23").Value = 1.02929036713324
$ws.Range("N23").Value = 1.03334578560852

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026295185622253
$ws.Range("D24").Value = 1.031941762114784
$ws.Range("E24").Value = 1.026614833363317
$ws.Range("I24").Value = 1.034970803324576
$ws.Range("J24").Value = 1.0325507262731
$ws.Range("K24").Value = 1.035326046794808
$ws.Range("L24").Value = 1.030017975971495
$ws.Range("N24").Value = 1.034017066739632

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027685104421995
$ws.Range("D25").Value = 1.032992345353342
$ws.Range("E25").Value = 1.027784584097642
$ws.Range("I25").Value = 1.035307314938674
$ws.Range("J25").Value = 1.033327211302193
$ws.Range("K25").Value = 1.036053724519098
$ws.Range("L25").Value = 1.030862518471634
$ws.Range("N25").Value = 1.03479465446653
